$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# --- New row of trade data (row 3) ---
# Plain text / time-like columns use the same style as row 2's A/B/C/D (s=2)
$ws.Range("A3").Value = "26-02-2026"
$ws.Range("B3").Value = "Thursday"
$ws.Range("C3").Value = "03:53"
$ws.Range("D3").Value = "17:36"
$ws.Range("A2:D2").Copy()
$ws.Range("A3:D3").PasteSpecial($xlPasteFormats)

# Numeric "$" columns use the same style as row 2's E..H (s=3)
$ws.Range("E3").Value = 68888.5
$ws.Range("F3").Value = 68800
$ws.Range("G3").Value = 70800
$ws.Range("H3").Value = 66600
$ws.Range("E2:H2").Copy()
$ws.Range("E3:H3").PasteSpecial($xlPasteFormats)

# CE/PE distance columns use style 2 (same as row 2's I/J)
$ws.Range("I3").Value = 10
$ws.Range("J3").Value = 11
$ws.Range("I2:J2").Copy()
$ws.Range("I3:J3").PasteSpecial($xlPasteFormats)

# Entry/exit premium columns use style 3 (same as row 2's K..P)
$ws.Range("K3").Value = 83
$ws.Range("L3").Value = 84
$ws.Range("M3").Value = 167
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 0
$ws.Range("K2:P2").Copy()
$ws.Range("K3:P3").PasteSpecial($xlPasteFormats)

# This trade is a WIN (positive P&L), so Q3/R3/S3 get a new green "profit"
# look instead of the red "loss" look used by row 2 — start from row 2's
# formatting (same number formats / border) then recolor font + fill.
$ws.Range("Q2").Copy()
$ws.Range("Q3").PasteSpecial($xlPasteFormats)
$ws.Range("Q3").Value = 167
$ws.Range("Q3").Font.Color = 24832        # font RGB 006100 (COM BGR order)
$ws.Range("Q3").Interior.Color = 13561798 # fill RGB C6EFCE (COM BGR order)

$ws.Range("R2").Copy()
$ws.Range("R3").PasteSpecial($xlPasteFormats)
$ws.Range("R3").Value = 15190.32
$ws.Range("R3").Font.Color = 24832
$ws.Range("R3").Interior.Color = 13561798

$ws.Range("S2").Copy()
$ws.Range("S3").PasteSpecial($xlPasteFormats)
$ws.Range("S3").Value = 100
$ws.Range("S3").Font.Color = 24832
$ws.Range("S3").Interior.Color = 13561798

$ws.Range("T3").Value = "Time Exit — Options Expired OTM (full premium kept)"
$ws.Range("U3").Value = "13h 43m"
$ws.Range("V3").Value = "DRY RUN"
$ws.Range("T2:V2").Copy()
$ws.Range("T3:V3").PasteSpecial($xlPasteFormats)

$ws.Range("W2").Copy()
$ws.Range("W3").PasteSpecial($xlPasteFormats)
$ws.Range("W3").Formula = "=W2+R3"

$excel.CutCopyMode = $false

Write-Output "row 3 added"
